# Apply the "Saldo" worksheet update described by the commit diff.
#
# The sheet ("Export") is a flat Conta/Nome/Saldo export with header on row 1
# and data rows starting at row 2. The edit:
#   - adds a new top balance row for account 004332783 (IRON) of 30009.06,
#     and removes the old, much-smaller row for the same account (9.06)
#     further down the list
#   - adds a new balance row for account 003553997 (MIRELLA) of 10448.06,
#     and removes the old row for the same account (448.06)
#   - adds a brand-new account 008071998 (ISADORA) with balance 100
#   - adds a new balance row for account 004643880 (GABRIEL) of 66.12,
#     and removes the old row for the same account (39.41)
#   - removes a stray duplicate-looking row for account 005616259 (MARIA, 31.25)
#
# All row numbers below are against the ORIGINAL (pre-edit) sheet, and the
# operations are applied from the bottom of the sheet upward so that each
# row number is still valid at the point it is used (deletions/insertions
# further down never shift the rows we still have to touch further up).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Conta numbers are zero-padded account numbers (e.g. "004332783"); they must
# be written as text or Excel's automatic type detection will strip the
# leading zeros and turn them into plain numbers. Force the Conta cell to
# Text format before assigning the value so the literal digit string is kept.
function Set-DataRow($row, $conta, $nome, $saldo) {
    $contaCell = $ws.Cells.Item($row, 1)
    $contaCell.NumberFormat = "@"
    $contaCell.Value = $conta
    $ws.Cells.Item($row, 2).Value = $nome
    $ws.Cells.Item($row, 3).Value = $saldo
}

# --- bottom-up ---

# Remove old 004332783 / IRON / 9.06 row
$ws.Rows(142).Delete()

# Remove stray 005616259 / MARIA / 31.25 row
$ws.Rows(103).Delete()

# Remove old 004643880 / GABRIEL / 39.41 row
$ws.Rows(88).Delete()

# Insert new 004643880 / GABRIEL / 66.12 row before row 60 (005558076)
$ws.Rows(60).Insert()
Set-DataRow 60 "004643880" "GABRIEL" 66.12

# Insert new 008071998 / ISADORA / 100 row before row 28 (004451652)
$ws.Rows(28).Insert()
Set-DataRow 28 "008071998" "ISADORA" 100

# Remove old 003553997 / MIRELLA / 448.06 row
$ws.Rows(16).Delete()

# Insert new 003553997 / MIRELLA / 10448.06 row before row 5 (004313254)
$ws.Rows(5).Insert()
Set-DataRow 5 "003553997" "MIRELLA" 10448.06

# Insert new 004332783 / IRON / 30009.06 row before row 4 (004368468)
$ws.Rows(4).Insert()
Set-DataRow 4 "004332783" "IRON" 30009.06
